# time_log.xlsx — "Finish redo & anki for clean up the words"
#
# Week 5 / JS101 log: the hours logged against 2021-10-22 (row 59, the day
# noted "Finished 2 small problems") are corrected from 1.5 to 2 hours.
# The weekly-total and grand-total SUBTOTAL formulas recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the logged hours for 2021-10-22 (Table1[Hours], row 59).
$ws.Range("C59").Value2 = 2

# Match the author's final cursor position/selection on save.
$ws.Range("D59").Select()
